# Auto-generated edit script: Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 27 de Marzo de 2020 a las 12:12'
$ws.Range("B7").Value = 64059
$ws.Range("C7").Value = 6273
$ws.Range("D7").Value = 9357
$ws.Range("E7").Value = 49844
$ws.Range("F7").Value = 4165
$ws.Range("G7").Value = 493
$ws.Range("H7").Value = 4858
$ws.Range("D11").Value = 897
$ws.Range("E11").Value = 10857
$ws.Range("B20").Value = 3441
$ws.Range("C20").Value = 69
$ws.Range("E20").Value = 3420
$ws.Range("B21").Value = 3180
$ws.Range("C21").Value = 130
$ws.Range("E21").Value = 2997
$ws.Range("B23").Value = 2988
$ws.Range("C23").Value = 3
$ws.Range("E23").Value = 2905
$ws.Range("A33").Value = 'Rumania'
$ws.Range("B33").Value = 1292
$ws.Range("C33").Value = 263
$ws.Range("D33").Value = 115
$ws.Range("E33").Value = 1153
$ws.Range("F33").Value = 32
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 24
$ws.Range("B34").Value = 1252
$ws.Range("C34").Value = 51
$ws.Range("E34").Value = 1220
$ws.Range("A35").Value = 'Polonia'
$ws.Range("B35").Value = 1244
$ws.Range("C35").Value = 23
$ws.Range("D35").Value = 7
$ws.Range("E35").Value = 1221
$ws.Range("F35").Value = 3
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 16
$ws.Range("A36").Value = 'Tailandia'
$ws.Range("B36").Value = 1136
$ws.Range("C36").Value = 91
$ws.Range("D36").Value = 97
$ws.Range("E36").Value = 1034
$ws.Range("F36").Value = 11
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 5
$ws.Range("A37").Value = 'Indonesia'
$ws.Range("B37").Value = 1046
$ws.Range("C37").Value = 153
$ws.Range("D37").Value = 46
$ws.Range("E37").Value = 913
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 9
$ws.Range("H37").Value = 87
$ws.Range("B38").Value = 1038
$ws.Range("C38").Value = 80
$ws.Range("E38").Value = 1023
$ws.Range("A39").Value = 'Rusia'
$ws.Range("B39").Value = 1036
$ws.Range("C39").Value = 196
$ws.Range("D39").Value = 45
$ws.Range("E39").Value = 988
$ws.Range("F39").Value = 8
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 3
$ws.Range("A62").Value = 'Libano'
$ws.Range("B62").Value = 391
$ws.Range("C62").Value = 23
$ws.Range("D62").Value = 23
$ws.Range("E62").Value = 361
$ws.Range("F62").Value = 3
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 7
$ws.Range("A63").Value = 'Irak'
$ws.Range("B63").Value = 382
$ws.Range("D63").Value = 105
$ws.Range("E63").Value = 241
$ws.Range("F63").Value = 0
$ws.Range("H63").Value = 36
$ws.Range("A79").Value = 'Ucrania'
$ws.Range("B79").Value = 226
$ws.Range("C79").Value = 30
$ws.Range("D79").Value = 4
$ws.Range("E79").Value = 217
$ws.Range("F79").Value = 0
$ws.Range("H79").Value = 5
$ws.Range("A80").Value = 'Kuwait'
$ws.Range("B80").Value = 225
$ws.Range("C80").Value = 17
$ws.Range("D80").Value = 57
$ws.Range("E80").Value = 168
$ws.Range("F80").Value = 11
$ws.Range("H80").Value = 0
$ws.Range("A81").Value = 'Principado de Andorra'
$ws.Range("B81").Value = 224
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 1
$ws.Range("E81").Value = 220
$ws.Range("F81").Value = 6
$ws.Range("H81").Value = 3
$ws.Range("A85").Value = 'Albania'
$ws.Range("B85").Value = 186
$ws.Range("C85").Value = 12
$ws.Range("D85").Value = 17
$ws.Range("E85").Value = 161
$ws.Range("F85").Value = 3
$ws.Range("G85").Value = 2
$ws.Range("H85").Value = 8
$ws.Range("A86").Value = 'Moldavia'
$ws.Range("B86").Value = 177
$ws.Range("D86").Value = 2
$ws.Range("E86").Value = 173
$ws.Range("F86").Value = 28
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 2
$ws.Range("A97").Value = 'Senegal'
$ws.Range("B97").Value = 119
$ws.Range("C97").Value = 14
$ws.Range("E97").Value = 108
$ws.Range("F97").Value = 0
$ws.Range("A98").Value = 'Brunei'
$ws.Range("B98").Value = 115
$ws.Range("C98").Value = 1
$ws.Range("D98").Value = 11
$ws.Range("E98").Value = 104
$ws.Range("F98").Value = 1
$ws.Range("H98").Value = 0
$ws.Range("A99").Value = 'Venezuela'
$ws.Range("B99").Value = 107
$ws.Range("D99").Value = 15
$ws.Range("E99").Value = 91
$ws.Range("F99").Value = 2
$ws.Range("H99").Value = 1
$ws.Range("A100").Value = 'Sri Lanka'
$ws.Range("B100").Value = 106
$ws.Range("D100").Value = 7
$ws.Range("E100").Value = 99
$ws.Range("F100").Value = 5
$ws.Range("A136").Value = 'Togo'
$ws.Range("C136").Value = 1
$ws.Range("D136").Value = 1
$ws.Range("E136").Value = 24
$ws.Range("H136").Value = 0
$ws.Range("A137").Value = 'Guatemala'
$ws.Range("B137").Value = 25
$ws.Range("D137").Value = 4
$ws.Range("E137").Value = 20
$ws.Range("H137").Value = 1
$ws.Range("A138").Value = 'Barbados'
$ws.Range("D138").Value = 0
$ws.Range("E138").Value = 24
